$d = $word.ActiveDocument

# Change 1: "delle quantità per ogni copia, e ogni info su di essa"
#        -> "delle quantità per ogni copia, autore, anno e descrizione"
$d.Content.Find.Execute("delle quantità per ogni copia, e ogni info su di essa", $true, $false, $false, $false, $false,
                         $true, 1, $false, "delle quantità per ogni copia, autore, anno e descrizione", 2)

# Change 2: "HTML/CSS – UI" -> "HTML/CSS – GUI"
$d.Content.Find.Execute("HTML/CSS – UI", $true, $false, $false, $false, $false,
                         $true, 1, $false, "HTML/CSS – GUI", 2)

# Change 3: add a new empty ListParagraph-styled paragraph after "Version’s System"
$r = $d.Content.Find.Execute("Version’s System", $true, $false, $false, $false, $false,
                              $true, 1, $false, "", 0)
$para = $d.Paragraphs.Last
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Version") {
        $p.Range.InsertParagraphAfter()
        $newPara = $d.Paragraphs.Item($i + 1)
        $newPara.Style = "List Paragraph"
        break
    }
}
